# Applies the "updated summary charts and summary reports including
# comments from Prof. Erhardt" revision to the FAC summary report.
#
# Sheet1 ("Sheet1"):
#   - Year 1 label (C1) and the mirrored header (E7) move from 2003 to 2012.
#   - The "Average Values" (E/F) and "Ridership Effect" (H) columns for every
#     factor row (8-21) switch to a #,##0.00 number format, and get refreshed
#     source data for the new 2012 base year.
#   - The "% Diff" columns (G/I) switch from a *100 "percent-as-number" style
#     formula to a true percentage (0.00%) format with a plain ratio formula.
#   - The selection / scroll position of the sheet view is updated.
#
# Sheet2 keeps its own "Average Values" header cell untouched in terms of
# content (only a cosmetic style-table reshuffle happens there, which the
# engine manages on its own when we touch styles elsewhere).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Header year: 2003 -> 2012 (kept as text, matching the original inline
# string cells)
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "'2012"
$ws1.Range("E7").Value = "'2012"

# ---------------------------------------------------------------------------
# Factor rows 8-18: refresh "Average Values" (E), keep F as-is, recompute the
# "Riddership Effect" (H), switch formulas in G/I to plain ratios, and apply
# the new number formats across E:I.
# ---------------------------------------------------------------------------
$rows = @(
  @{ Row = 8;  E = 800855;       H = -712095.4942600001 },
  @{ Row = 9;  E = 0.355993915;  H = -926194.1929189999 },
  @{ Row = 10; E = 355244.33;    H = 21898.571044 },
  @{ Row = 11; E = 19.93181551;  H = 1592.691603100001 },
  @{ Row = 12; E = 4.3491;       H = -39594.271964 },
  @{ Row = 13; E = 22210.93;     H = -471.1802779999998 },
  @{ Row = 14; E = 4.88;         H = -18228.108403 },
  @{ Row = 15; E = 3.9;          H = -1012.612951500001 },
  @{ Row = 16; E = $null;        H = -87509.9669 },
  @{ Row = 17; E = 0;            H = 0 },
  @{ Row = 18; E = 0;            H = 0 }
)

foreach ($r in $rows) {
  $row = $r.Row

  if ($null -ne $r.E) {
    $ws1.Range("E$row").Value = $r.E
  }
  $ws1.Range("H$row").Value = $r.H

  $ws1.Range("G$row").Formula = "=IFERROR((F$row-E$row)/E$row,0)"
  $ws1.Range("I$row").Formula = "=IFERROR(H$row/`$E`$21,0)"

  $ws1.Range("E$row:F$row").NumberFormat = "#,##0.00"
  $ws1.Range("H$row").NumberFormat = "#,##0.00"
  $ws1.Range("G$row").NumberFormat = "0.00%"
  $ws1.Range("I$row").NumberFormat = "0.00%"
}

# ---------------------------------------------------------------------------
# Row 19 "New Reporters": gains an explicit 0 in H19 and the same number
# format treatment as the rows above.
# ---------------------------------------------------------------------------
$ws1.Range("H19").Value = 0
$ws1.Range("G19").Formula = "=IFERROR((F19-E19)/E19,0)"
$ws1.Range("I19").Formula = "=IFERROR(H19/`$E`$21,0)"
$ws1.Range("E19:F19").NumberFormat = "#,##0.00"
$ws1.Range("H19").NumberFormat = "#,##0.00"
$ws1.Range("G19").NumberFormat = "0.00%"
$ws1.Range("I19").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Row 20 "Total Modeled Ridership": refreshed base-year total + new formula
# style for G20 (I20 stays "=G20").
# ---------------------------------------------------------------------------
$ws1.Range("E20").Value = 1289738.196
$ws1.Range("G20").Formula = "=IFERROR((F20-E20)/E20,0)"
$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat = "#,##0.00"
$ws1.Range("G20").NumberFormat = "0.00%"
$ws1.Range("I20").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Row 21 "Total Observed Ridership": refreshed base-year total + new formula
# style for G21 (I21 stays "=G21").
# ---------------------------------------------------------------------------
$ws1.Range("E21").Value = 1538794
$ws1.Range("G21").Formula = "=IFERROR((F21-E21)/E21,0)"
$ws1.Range("E21:F21").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat = "#,##0.00"
$ws1.Range("G21").NumberFormat = "0.00%"
$ws1.Range("I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Sheet view: drop the frozen scroll position (topLeftCell) and move the
# active selection to H21.
# ---------------------------------------------------------------------------
$ws1.Range("H21").Select()
